$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "What is AI" paragraph (paragraph 10): merge the " - " run and the
#    "The ability of machine..." run into a single run, producing
#    " - The ability of machine to perform tasks that require human
#    intelligence " while leaving the leading "What is AI" run untouched.
# ---------------------------------------------------------------------------
$whatIsAiPara = $d.Paragraphs(10)
$null = $whatIsAiPara.Range.Find.Execute(
    " `u{2013} The ability of machine to perform tasks that require human intelligence ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " `u{2013} The ability of machine to perform tasks that require human intelligence ",
    2)

# ---------------------------------------------------------------------------
# 2) Insert two new, empty "ListParagraph" paragraphs right after it (mirrors
#    the paragraph that used to hold only the underline run).
# ---------------------------------------------------------------------------
$whatIsAiPara.Range.InsertParagraphAfter()
$blank1 = $d.Paragraphs(11)
$blank1.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:spacing w:lineRule='auto' w:line='240' w:before='0' w:after='0'/><w:ind w:hanging='360' w:left='709' w:right='0'/><w:contextualSpacing/><w:rPr><w:color w:themeColor='text1' w:val='000000'/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p>")

$blank1.Range.InsertParagraphAfter()
$blank2 = $d.Paragraphs(12)
$blank2.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:spacing w:lineRule='auto' w:line='240' w:before='0' w:after='0'/><w:ind w:hanging='0' w:left='709' w:right='0'/><w:contextualSpacing/></w:pPr><w:r><w:rPr/></w:r></w:p>")

# ---------------------------------------------------------------------------
# 3) Replace the old single-run "underline" paragraph (now paragraph 13)
#    with the new "AI Algorithms ..." bullet line, restyled to match the
#    other ListParagraph bullets.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs(13)
$lastPara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:spacing w:lineRule='auto' w:line='240' w:before='0' w:after='0'/><w:ind w:hanging='360' w:left='709' w:right='0'/><w:contextualSpacing/><w:rPr><w:color w:themeColor='text1' w:val='000000'/></w:rPr></w:pPr><w:r><w:rPr><w:color w:themeColor='text1' w:val='000000'/></w:rPr><w:t xml:space='preserve'>AI Algorithms can recognize speech, understand natural language, make decisions, Learn from data </w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 4) Styles: duplicate "Footnote Characters1"/"Endnote Characters1" into new
#    "...2" character styles (superscript, quick-style).
# ---------------------------------------------------------------------------
$footChar2 = $d.Styles.Add("FootnoteCharacters2", 2)
$footChar2.NameLocal = "Footnote Characters2"
$footChar2.QuickStyle = $true
$footChar2.Font.Superscript = $true

$endChar2 = $d.Styles.Add("EndnoteCharacters2", 2)
$endChar2.NameLocal = "Endnote Characters2"
$endChar2.QuickStyle = $true
$endChar2.Font.Superscript = $true

Write-Output "done"
